# Edit Sheet1 in place to keep only the "Gamma-HCH" / "sum_DDT" summary rows
# (dropping ESTUARY), then add a new Sheet2 that holds the original
# per-estuary breakdown with updated "First/Last 5 years" statistics.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Sheet2: per-estuary data (ESTUARY, PARAMETRE_LIBELLE, First 5 years, Last 5 years) ---
$ws2.Range("A1").Value = "ESTUARY"
$ws2.Range("B1").Value = "PARAMETRE_LIBELLE"
$ws2.Range("C1").Value = "First 5 years"
$ws2.Range("D1").Value = "Last 5 years"

$ws2.Range("A2").Value = "Gironde"
$ws2.Range("B2").Value = "Gamma-HCH"
$ws2.Range("C2").Value = 0.9743000000000001
$ws2.Range("D2").Value = 0

$ws2.Range("A3").Value = "Loire"
$ws2.Range("B3").Value = "Gamma-HCH"
$ws2.Range("C3").Value = 1.3579
$ws2.Range("D3").Value = 0

$ws2.Range("A4").Value = "Seine"
$ws2.Range("B4").Value = "Gamma-HCH"
$ws2.Range("C4").Value = 2.2463
$ws2.Range("D4").Value = 0

$ws2.Range("A5").Value = "Gironde"
$ws2.Range("B5").Value = "sum_DDT"
$ws2.Range("C5").Value = 15.77605
$ws2.Range("D5").Value = 1.737

$ws2.Range("A6").Value = "Loire"
$ws2.Range("B6").Value = "sum_DDT"
$ws2.Range("C6").Value = 12.28675
$ws2.Range("D6").Value = 0.5519000000000001

$ws2.Range("A7").Value = "Seine"
$ws2.Range("B7").Value = "sum_DDT"
$ws2.Range("C7").Value = 23.24085
$ws2.Range("D7").Value = 1.7117

# Header style (bold, centered) matching Sheet1's header
$ws2.Range("A1:D1").Font.Bold = $true
$ws2.Range("A1:D1").HorizontalAlignment = -4108

# --- Sheet1: clear everything, then write the new summarised table ---
$ws1.Cells.Clear()

$ws1.Range("A1").Value = "PARAMETRE_LIBELLE"
$ws1.Range("B1").Value = "First 5 years"
$ws1.Range("C1").Value = "Last 5 years"

$ws1.Range("A2").Value = "Gamma-HCH"
$ws1.Range("B2").Value = 1.3865
$ws1.Range("C2").Value = 0

$ws1.Range("A3").Value = "sum_DDT"
$ws1.Range("B3").Value = 15.63015
$ws1.Range("C3").Value = 1.112

$ws1.Range("A1:C1").Font.Bold = $true
$ws1.Range("A1:C1").HorizontalAlignment = -4108

# Keep Sheet1 as the active/selected tab (matches original workbook state)
$ws1.Activate()
